$wb = $excel.ActiveWorkbook

# Mapping of row -> new F-column value, applied to both the "展览" and
# "全部类型" worksheets (they contain duplicated data tables).
$updates = @{
    2  = 1529
    4  = 989
    6  = 2437
    8  = 1502
    12 = 440
    14 = 19
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
